# Fix voor relaties Excelbestanden
# Adds 5 new "Concepten" relation rows (focus, hetzelfde, gerelateerd, brederdan,
# engerdan) to the "URI schema" sheet, right after the existing "betrekkingop"
# row, and repairs / re-creates all hyperlinks in column D so they keep
# pointing at the correct (shifted) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URI schema")

# ---------------------------------------------------------------------------
# 1. Insert 5 new rows right before the current row 17 ("Collecties" block).
#    This shifts the old rows 17-26 down to 22-31 and extends the bottom
#    blank filler rows from 44 down to 49 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("17:21").Insert()

# ---------------------------------------------------------------------------
# 2. Populate the newly inserted rows with the new "Concepten" relations.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 17; B = "focus";       D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{focus}" },
    @{ Row = 18; B = "hetzelfde";   D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{hetzelfde}" },
    @{ Row = 19; B = "gerelateerd"; D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{gerelateerd}" },
    @{ Row = 20; B = "brederdan";   D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{brederdan}" },
    @{ Row = 21; B = "engerdan";    D = "http://data.test.pdok.nl/catalogus/dso/id/concept/{engerdan}" }
)

foreach ($r in $newRows) {
    $ws.Range("A" + $r.Row).Value = "Concepten"
    $ws.Range("B" + $r.Row).Value = $r.B
    $ws.Range("D" + $r.Row).Value = $r.D
}

# ---------------------------------------------------------------------------
# 3. Rebuild the hyperlinks collection. Inserting rows does not move the
#    existing Hyperlink objects along with their cells in this engine, so we
#    drop the whole collection and recreate every hyperlink at its correct
#    (possibly shifted) address, preserving address/location/display exactly
#    as before where it did not move.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ref, address, subaddress(location), displayText (display text that differs
# from the cell's own value - only set when the original hyperlink carried a
# "display" attribute different from the stored cell text)
$links = @(
    @{ Ref = "D5";  Address = "http://purl.org/dc/terms/{eigenschap}" },
    @{ Ref = "D6";  Address = "http://purl.org/iso25964/skos-thes"; SubAddress = "{eigenschap}" },
    @{ Ref = "D12"; Address = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{specialisatie}" },
    @{ Ref = "D13"; Address = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{generalisatie}" },
    @{ Ref = "D14"; Address = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{onderdeel}" },
    @{ Ref = "D15"; Address = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{bestaatuit}" },
    @{ Ref = "D16"; Address = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/catalogus/dsoprogramma/id/begrip/{betrekkingop}" },
    @{ Ref = "D22"; Address = "http://localhost:8080/catalogus/dso/id/collection/"; Display = "http://localhost:8080/catalogus/dso/id/collection/{collectie}" },
    @{ Ref = "D23"; Address = "http://localhost:8080/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/dsoprogramma/id/begrip/{begrip}" },
    @{ Ref = "D9";  Address = "http://localhost:8080/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/dsoprogramma/id/begrip/{begrip}" },
    @{ Ref = "D7";  Address = "http://www.w3.org/1999/02/22-rdf-syntax-ns"; SubAddress = "{eigenschap}" },
    @{ Ref = "D2";  Address = "http://www.w3.org/2004/02/skos/core"; SubAddress = "{klasse}" },
    @{ Ref = "D11"; Address = "http://localhost:8080/catalogus/dso/id/concept/"; Display = "http://localhost:8080/catalogus/dso/id/concept/{bron}" },
    @{ Ref = "D28"; Address = "http://purl.org/dc/dcmitype/{subklasse}" },
    @{ Ref = "D27"; Address = "http://localhost:8080/catalogus/dso/id/concept/"; Display = "http://localhost:8080/catalogus/dso/id/concept/{bron}" },
    @{ Ref = "D25"; Address = "http://localhost:8080/catalogus/dso/id/collection/"; Display = "http://localhost:8080/catalogus/dso/id/collection/{waardelijst}" },
    @{ Ref = "D26"; Address = "http://localhost:8080/catalogus/dso/id/collection/"; Display = "http://localhost:8080/catalogus/dso/id/collection/{collectie}" },
    @{ Ref = "D30"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{toeleidingsbegrip}" },
    @{ Ref = "D31"; Address = "http://localhost:8080/dsoprogramma/id/begrip/"; Display = "http://localhost:8080/dsoprogramma/id/begrip/{begrip}" },
    @{ Ref = "D8";  Address = "http://xmlns.com/foaf/0.1/{eigenschap}" },
    @{ Ref = "D10"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{domein}" },
    @{ Ref = "D24"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{domein}" },
    @{ Ref = "D17"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{focus}" },
    @{ Ref = "D18"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{hetzelfde}" },
    @{ Ref = "D19"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{gerelateerd}" },
    @{ Ref = "D20"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{brederdan}" },
    @{ Ref = "D21"; Address = "http://data.test.pdok.nl/catalogus/dso/id/concept/{engerdan}" }
)

foreach ($lnk in $links) {
    $cell = $ws.Range($lnk.Ref)
    $originalValue = $cell.Value()

    $sub = $null
    if ($lnk.ContainsKey("SubAddress")) { $sub = $lnk.SubAddress }

    if ($lnk.ContainsKey("Display")) {
        $ws.Hyperlinks.Add($cell, $lnk.Address, $sub, $null, $lnk.Display)
        # Adding a hyperlink with display text overwrites the cell's value,
        # so restore the original (possibly different) stored text.
        $cell.Value = $originalValue
    } else {
        $ws.Hyperlinks.Add($cell, $lnk.Address, $sub)
    }

    # Adding a hyperlink re-applies the "Hyperlink" cell style but does so
    # through a freshly duplicated style record; explicitly re-assigning the
    # named style collapses it back onto the workbook's existing "Hyperlink"
    # style (matching the original file's cell formatting exactly).
    $cell.Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 4. Restore the active cell/selection shown in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
